$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new favourite/shopping "count" value in F1
$ws.Range("F1").Value = 58

# Move the active selection to G1 (next empty cell), mirroring the
# after-entry cursor position recorded in the saved workbook
$ws.Range("G1").Select()
